$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 5
$ws.Range("Q2").Value = 3.4
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05
$ws.Range("N4").Value = 8
$ws.Range("O7").Value = 1.17
$ws.Range("P7").Value = 5
$ws.Range("Q7").Value = 1.6
$ws.Range("R7").Value = 2.3
$ws.Range("W7").Value = 11
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 11
$ws.Range("AK7").Value = 21
$ws.Range("AQ7").Value = 41
$ws.Range("AX7").Value = 15
$ws.Range("G8").Value = 2.15
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 2.75
$ws.Range("L8").Value = 3.6
$ws.Range("Q8").Value = 1.73
$ws.Range("R8").Value = 2.08
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.2
$ws.Range("W8").Value = 9.5
$ws.Range("Y8").Value = 9
$ws.Range("AE8").Value = 13
$ws.Range("AJ8").Value = 34
$ws.Range("AK8").Value = 23
$ws.Range("AL8").Value = 29
$ws.Range("AN8").Value = 4.33
$ws.Range("AO8").Value = 11
$ws.Range("AV8").Value = 51
$ws.Range("AW8").Value = 5.5
$ws.Range("AX8").Value = 17
$ws.Range("AY8").Value = 23
$ws.Range("BB8").Value = 151
$ws.Range("G13").Value = 1.65
$ws.Range("I13").Value = 6
$ws.Range("L13").Value = 7
$ws.Range("AI13").Value = 21
$ws.Range("AU13").Value = 11
$ws.Range("G17").Value = 3.2
$ws.Range("I17").Value = 2.3
$ws.Range("K17").Value = 2.05
$ws.Range("Y17").Value = 12
$ws.Range("AA17").Value = 29
$ws.Range("AC17").Value = 8
$ws.Range("AD17").Value = 6
$ws.Range("AH17").Value = 10
$ws.Range("AI17").Value = 9.5
$ws.Range("AJ17").Value = 21
$ws.Range("AO17").Value = 19
$ws.Range("AZ17").Value = 41
$ws.Range("M20").Value = 1.05
